$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1 / K1 were shared-string text ("r" / "s"); retrain replaced them with
# actual numeric values consistent with the rest of the column.
$ws.Range("J1").Value = 0.3
$ws.Range("K1").Value = 1

# Column K (rows 2-51) retrained from 0.3 to 1. Column J is left untouched.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 11).Value = 1
}

# Update the view state to match what was saved (scrolled one row further,
# selection moved to K1:K51 with active cell K1).
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("K1:K51").Select()
